$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect it so we can update cell values,
# then re-apply protection afterwards with the same settings.
$ws.Unprotect()

# Update the confidential/disclosure banner date (2021-05-07 -> 2021-05-10)
$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-10 for illustrative purposes only and are subject to change."

# Update the Weight (col D) and Percent Change (col E) figures for each holding
$ws.Range("D2").Value = 0.06367875953186783
$ws.Range("E2").Value = -0.02580446970278805
$ws.Range("D3").Value = 0.03832739817048127
$ws.Range("E3").Value = -0.03072052885973742
$ws.Range("D4").Value = 0.03233599495788284
$ws.Range("E4").Value = -0.02091420423037316
$ws.Range("D5").Value = 0.02973077378445601
$ws.Range("E5").Value = -0.02649491642254009
$ws.Range("D6").Value = 0.0273857952731642
$ws.Range("E6").Value = -0.02558749622650325
$ws.Range("D7").Value = 0.02503297649076856
$ws.Range("E7").Value = -0.0001240387000744914
$ws.Range("D8").Value = 0.1879142205408955
$ws.Range("E8").Value = -0.002375296912114022
$ws.Range("D9").Value = 0.02419810405787902
$ws.Range("E9").Value = 0.01050445103857567
$ws.Range("D10").Value = 0.0220315499340346
$ws.Range("E10").Value = 0.01864594894561589
$ws.Range("D11").Value = 0.02210940925994703
$ws.Range("E11").Value = -0.00907246030798603
$ws.Range("D12").Value = 0.02008786133830043
$ws.Range("E12").Value = -0.002921445574550985
$ws.Range("D13").Value = 0.01997312945026596
$ws.Range("E13").Value = -0.003082029397818986
$ws.Range("D14").Value = 0.01721246131757728
$ws.Range("E14").Value = -0.005053340819764163
$ws.Range("D15").Value = 0.01613853824246708
$ws.Range("E15").Value = -0.02274170274170273
$ws.Range("D16").Value = 0.01512784190809035
$ws.Range("E16").Value = -0.04433497536945807
$ws.Range("D17").Value = 0.01414819615234271
$ws.Range("E17").Value = 0.004422253922967245
$ws.Range("D18").Value = 0.01412413195390517
$ws.Range("E18").Value = 0.007969222313822577
$ws.Range("D19").Value = 0.01362297561483173
$ws.Range("E19").Value = -0.0410868747649491
$ws.Range("D20").Value = 0.01356941336669656
$ws.Range("E20").Value = 0.002402691013935643
$ws.Range("D21").Value = 0.01235750928280336
$ws.Range("E21").Value = 0.01461442786069678
$ws.Range("D22").Value = 0.01302862872664775
$ws.Range("E22").Value = 0.004060475161987043
$ws.Range("D23").Value = 0.01165638721735861
$ws.Range("E23").Value = -0.02413425679275438
$ws.Range("D24").Value = 0.01296827416443747
$ws.Range("E24").Value = 0.007757691847240622
$ws.Range("D25").Value = 0.01185064739990683
$ws.Range("E25").Value = 0.005512159174649867
$ws.Range("D26").Value = 0.009344050626850633
$ws.Range("E26").Value = -0.059665038381019
$ws.Range("D27").Value = 0.009789315924391716
$ws.Range("E27").Value = -0.06992419196244481
$ws.Range("D28").Value = 0.01035536797279994
$ws.Range("E28").Value = -0.01652173913043475
$ws.Range("D29").Value = 0.00987823701894076
$ws.Range("E29").Value = -0.01790127579987977
$ws.Range("D30").Value = 0.009847691012214397
$ws.Range("E30").Value = -0.0127699826580483
$ws.Range("D31").Value = 0.008729598489004323
$ws.Range("E31").Value = -0.02947806485174276
$ws.Range("D32").Value = 0.01037500746378284
$ws.Range("E32").Value = -0.01032524522457401
$ws.Range("D33").Value = 0.009344361132636925
$ws.Range("E33").Value = 0.01362397820163497
$ws.Range("D34").Value = 0.008885976965625054
$ws.Range("E34").Value = 0.007338103100348548
$ws.Range("D35").Value = 0.009248414844673052
$ws.Range("E35").Value = 0.003273459795198885
$ws.Range("D36").Value = 0.008310299237841812
$ws.Range("E36").Value = -0.00008406893652790348
$ws.Range("D37").Value = 0.00858859004880499
$ws.Range("E37").Value = -0.0006778741865511861
$ws.Range("D38").Value = 0.007829054082314277
$ws.Range("E38").Value = -0.06444368428097635
$ws.Range("D39").Value = 0.008602873314974368
$ws.Range("E39").Value = 0.00707427993936327
$ws.Range("D40").Value = 0.008113438569333413
$ws.Range("E40").Value = -0.004090165424468384
$ws.Range("D41").Value = 0.007026435437975708
$ws.Range("E41").Value = -0.03641345176543365
$ws.Range("D42").Value = 0.00749056396203396
$ws.Range("E42").Value = -0.06463547334058761
$ws.Range("D43").Value = 0.008113516195779984
$ws.Range("E43").Value = -0.003109452736318463
$ws.Range("D44").Value = 0.007488390421529925
$ws.Range("E44").Value = -0.01023147812205205
$ws.Range("D45").Value = 0.0072728217793975
$ws.Range("E45").Value = -0.01063080371437708
$ws.Range("D46").Value = 0.007889952029650571
$ws.Range("E46").Value = -0.002597402597402709
$ws.Range("D47").Value = 0.007304027610919731
$ws.Range("E47").Value = -0.003060834077286079
$ws.Range("D48").Value = 0.007200551557638303
$ws.Range("E48").Value = -0.01016612943218431
$ws.Range("D49").Value = 0.006656855925843022
$ws.Range("E49").Value = -0.0003498338289313896
$ws.Range("D50").Value = 0.007267465554583983
$ws.Range("E50").Value = 0.008694630478204513
$ws.Range("D51").Value = 0.006547519075845354
$ws.Range("E51").Value = 0.006722270633636285
$ws.Range("D52").Value = 0.006783231780863398
$ws.Range("E52").Value = -0.01874506482954352
$ws.Range("D53").Value = 0.005369033177201709
$ws.Range("E53").Value = -0.01474734330947725
$ws.Range("D54").Value = 0.006067904075695845
$ws.Range("E54").Value = -0.00337734111145227
$ws.Range("D55").Value = 0.006047294254130791
$ws.Range("E55").Value = 0.005583902955617503
$ws.Range("D56").Value = 0.0056810216286221
$ws.Range("E56").Value = -0.0009797210049723448
$ws.Range("D57").Value = 0.006753500851826049
$ws.Range("E57").Value = -0.006298850574712689
$ws.Range("D58").Value = 0.005449376549404467
$ws.Range("E58").Value = 0.00786324786324788
$ws.Range("D59").Value = 0.005245374247811377
$ws.Range("E59").Value = 0.0006659563132658253
$ws.Range("D60").Value = 0.00484947937029054
$ws.Range("E60").Value = 0.02253809706748622
$ws.Range("D61").Value = 0.004869739872846017
$ws.Range("E61").Value = -0.005435735577766221
$ws.Range("D62").Value = 0.004850488514095985
$ws.Range("E62").Value = -0.0007201728414820696
$ws.Range("D63").Value = 0.004207353404241056
$ws.Range("E63").Value = 0.01033210332103329
$ws.Range("D64").Value = 0.004001798573716511
$ws.Range("E64").Value = 0.004034761018001376
$ws.Range("D65").Value = 0.003716599009008363
$ws.Range("E65").Value = 0.03024353565311833
$ws.Range("D66").Value = 0.003671109911316753
$ws.Range("E66").Value = 0
$ws.Range("D67").Value = 0.003805481290334119
$ws.Range("E67").Value = 0.001529894131326204
$ws.Range("D68").Value = 0.003664239970795067
$ws.Range("E68").Value = -0.01722329911976861
$ws.Range("D69").Value = 0.003610250777203746
$ws.Range("E69").Value = 0.01051432011696929
$ws.Range("D70").Value = 0.003002280447646256
$ws.Range("E70").Value = -0.0267349260523323
$ws.Range("D71").Value = 0.002950115475549393
$ws.Range("E71").Value = -0.03840385222608156
$ws.Range("D72").Value = 0.002299644666493342
$ws.Range("E72").Value = -0.03689513747067463
$ws.Range("D73").Value = 0.001955565442059791
$ws.Range("E73").Value = -0.0340385836773579
$ws.Range("D74").Value = 0.001896918661674106
$ws.Range("E74").Value = -0.01913121764573489
$ws.Range("D75").Value = 0.001454098597200064
$ws.Range("E75").Value = -0.01831091180866984
$ws.Range("D76").Value = 0.001686745057578485
$ws.Range("E76").Value = -0.02324083022688561
$ws.Range("E77").Value = -0.01013022798980512

# Re-protect the worksheet (restores protection, matching the workbook's
# original protected state).
$ws.Protect()
